$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a row at 12 so the old "..."/"Task2.n" rows (12,13) shift down
# to (13,14), leaving room for the new "do the view" / "do the model"
# sub-tasks that get attached to rows 10 and 11.
$ws.Rows.Item(12).Insert()

# --- Register the custom date style (numFmtId 164) first, before any of
# the quote-prefixed text cells below, so it lands at cellXfs index 4 -
# matching the style used later by C22.
$ws.Range("C22").NumberFormat = "dd/MM/yyyy"

# --- New sub-task labels on the existing Goal2 rows.
$ws.Range("B10").Value = "do the view"
$ws.Range("B11").Value = "do the model"

# --- New "Goal" form-submission log rows (16, 18, 20, 22), mirroring a
# little append-only log of goal entries. Rows 15/17/19/21 stay empty.
$ws.Range("A16").Value = "Goal"
$ws.Range("B16").Value = "description"
$ws.Range("C16").Value = 44236.44610025463
$ws.Range("D16").Value = "00:00"
$ws.Range("E16").Value = "00:00"
$cell = $ws.Range("F16")
$cell.Value = "'0%"
$cell.Style = "Normal"
$ws.Range("G16").Value = 0.0

$ws.Range("A18").Value = "Goal"
$ws.Range("B18").Value = "description"
$ws.Range("C18").Value = 44245.44610025463
$ws.Range("D18").Value = "00:00"
$ws.Range("E18").Value = "00:00"
$cell = $ws.Range("F18")
$cell.Value = "'0%"
$cell.Style = "Normal"
$ws.Range("G18").Value = 0.0

$ws.Range("A20").Value = "Goal"
$ws.Range("B20").Value = "a"
$ws.Range("C20").Value = 44237.448924537035
$ws.Range("D20").Value = "00:00"
$ws.Range("E20").Value = "00:00"
$cell = $ws.Range("F20")
$cell.Value = "'0%"
$cell.Style = "Normal"
$ws.Range("G20").Value = 0.0

$ws.Range("A22").Value = "Goal"
$ws.Range("B22").Value = "a"
$ws.Range("C22").Value = 44236.456396574074
$ws.Range("D22").Value = "00:00"
$ws.Range("E22").Value = "00:00"
$cell = $ws.Range("F22")
$cell.Value = "'0%"
$cell.Style = "Normal"
$ws.Range("G22").Value = 0.0

# --- Column A is new (the Goal log sits outside the old B:G task table);
# give it the same "customWidth" treatment as the rest of the sheet.
$ws.Columns.Item(1).ColumnWidth = 11.3

# --- Leave the selection where the user ended up after adding the log.
$ws.Range("G3").Select()
